# "Added monte carlo emulation" - append 5 new task rows (TASK-7..TASK-11)
# to the "tasks" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")

# Row 14 - TASK-7 / "Some story just in case " / blocked by TASK-6
$ws.Range("A14").Value = "TASK-7"
$ws.Range("B14").Value = "Some story just in case "
$ws.Range("C14").Value = "TASK-6"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 13

# Row 15 - TASK-8 / "Tests for some story just in case " / blocked by TASK-7, TASK-2
$ws.Range("A15").Value = "TASK-8"
$ws.Range("B15").Value = "Tests for some story just in case "
$ws.Range("C15").Value = "TASK-7, TASK-2"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 5

# Row 16 - TASK-9 / blocked by TASK-1
$ws.Range("A16").Value = "TASK-9"
$ws.Range("C16").Value = "TASK-1"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8

# Row 17 - TASK-10 / blocked by TASK-9 (fill B17 last, after row 18, to
# preserve the original shared-string insertion order)
$ws.Range("A17").Value = "TASK-10"
$ws.Range("C17").Value = "TASK-9"
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 20

# Row 18 - TASK-11
$ws.Range("A18").Value = "TASK-11"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 3

# "Add new page" comment for TASK-10, filled last so it becomes the final
# new unique shared string added to the workbook.
$ws.Range("B17").Value = "Add new page"

# Update the sheet's active selection to match the authored state
$ws.Range("B16").Select() | Out-Null
